$d = $word.ActiveDocument

# Merge the "Versi" + "on" runs into a single "Version" run.
# (Word only merges runs across a Range.Text assignment when the text
# actually changes, so nudge through a distinct value first.)
$rVersion = $d.Range(0, 7)
$rVersion.Text = "Versionx"
$rVersion2 = $d.Range(0, 8)
$rVersion2.Text = "Version"

# " 2" -> " 1." (the run following the spellEnd proofErr mark)
$rNum = $d.Range(8, 9)
$rNum.Text = "1."

# Remove the trailing "." run entirely (now at offset 10..11)
$rDot = $d.Range(10, 11)
$rDot.Text = ""
